# Update RandomForest imputation results (re-run of the algorithm produced
# slightly different numeric outputs for a subset of cells in columns A, D
# and E). Apply the new values cell-by-cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.4454
$ws.Range("D7").Value = -7.2813
$ws.Range("A9").Value = -21.78519999999999
$ws.Range("D12").Value = -6.905600000000002
$ws.Range("A13").Value = -21.9919
$ws.Range("D14").Value = -7.751200000000007
$ws.Range("E15").Value = 16.24159999999999
$ws.Range("A16").Value = -21.46939999999998
$ws.Range("A18").Value = -22.21059999999999
$ws.Range("D19").Value = -7.894299999999999
$ws.Range("A20").Value = -19.62139999999998
$ws.Range("A26").Value = -20.99289999999996
$ws.Range("D26").Value = -8.407099999999996
$ws.Range("A27").Value = -21.56869999999998
$ws.Range("D27").Value = -8.637599999999992
$ws.Range("E28").Value = 16.7223
$ws.Range("A29").Value = -21.02729999999998
$ws.Range("D29").Value = -7.426500000000005
$ws.Range("E33").Value = 17.01150000000002
$ws.Range("A35").Value = -19.48519999999999
$ws.Range("E35").Value = 16.52440000000001
$ws.Range("A36").Value = -20.0665
$ws.Range("D37").Value = -8.078499999999998
$ws.Range("D38").Value = -7.487799999999996
$ws.Range("E38").Value = 16.89269999999999
$ws.Range("E43").Value = 17.16450000000001
$ws.Range("E44").Value = 16.26159999999998
$ws.Range("A45").Value = -21.65279999999998
$ws.Range("E45").Value = 16.7366
$ws.Range("D47").Value = -7.897999999999999
$ws.Range("E47").Value = 16.80470000000001
$ws.Range("D51").Value = -8.421699999999998
$ws.Range("E51").Value = 16.3506
$ws.Range("D52").Value = -7.707499999999999
$ws.Range("E54").Value = 16.6423
$ws.Range("A55").Value = -22.01529999999999
$ws.Range("D55").Value = -8.2781
$ws.Range("A57").Value = -22.2551
$ws.Range("E57").Value = 16.4555
$ws.Range("E62").Value = 16.307
$ws.Range("E63").Value = 18.20380000000002
$ws.Range("E67").Value = 17.10270000000002
$ws.Range("A69").Value = -21.60879999999999
$ws.Range("D69").Value = -7.284699999999994
$ws.Range("D70").Value = -7.569300000000003
$ws.Range("E70").Value = 17.1963
$ws.Range("A76").Value = -20.22899999999997
$ws.Range("D76").Value = -7.446499999999999
$ws.Range("A78").Value = -20.28279999999998
$ws.Range("D81").Value = -8.323199999999996
$ws.Range("E81").Value = 16.4327
$ws.Range("A82").Value = -22.0047
$ws.Range("A83").Value = -22.184
$ws.Range("D83").Value = -8.391300000000001
$ws.Range("E88").Value = 16.3471
$ws.Range("A93").Value = -21.07739999999998
$ws.Range("D94").Value = -7.088800000000004
$ws.Range("E96").Value = 15.1483
$ws.Range("A97").Value = -22.035
$ws.Range("E99").Value = 16.6505
$ws.Range("D100").Value = -8.3386
$ws.Range("D102").Value = -7.949899999999996
